$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three trailing balance rows (at the bottom of the data block,
# right before the blank row / "Filtros aplicados" footer row).
$ws.Rows.Item(232).Resize(3).Delete()

# Delete the ERIK row further up in the list.
$ws.Rows.Item(22).Delete()

# Insert a new row holding the EDINARDO record right above the KAUANNE row.
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "004211807"
$ws.Cells.Item(4, 1).ClearFormats()
$ws.Cells.Item(4, 2).Value = "EDINARDO"
$ws.Cells.Item(4, 3).Value = 22089.92
